$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H43").Value = 2895.0588
$ws.Range("I43").Value = 2796.6365
$ws.Range("K43").Value = 2796.6365
$ws.Range("M43").Value = -2727.6365
$ws.Range("H103").Value = 83334080
$ws.Range("J103").Value = 166667330
$ws.Range("L103").Value = 500001990
$ws.Range("N103").Value = -500003162
$ws.Range("H112").Value = 2651.4546
$ws.Range("J112").Value = 2929.889
$ws.Range("L112").Value = 8789.667000000001
$ws.Range("N112").Value = -11005.667
$ws.Range("H116").Value = 17906.285
$ws.Range("I116").Value = 10632
$ws.Range("J116").Value = 23362
$ws.Range("K116").Value = 10632
$ws.Range("L116").Value = 23362
$ws.Range("M116").Value = -7190
$ws.Range("N116").Value = -30246
$ws.Range("H132").Value = 3229.276
$ws.Range("I132").Value = 1987.037
$ws.Range("K132").Value = 5961.111
$ws.Range("M132").Value = -3431.111
$ws.Range("H138").Value = 2580.3333
$ws.Range("I138").Value = 1194.8182
$ws.Range("J138").Value = 5018.84
$ws.Range("K138").Value = 3584.4546
$ws.Range("L138").Value = 15056.52
$ws.Range("M138").Value = 1555.5454
$ws.Range("N138").Value = -25336.52

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 531.5
$ws.Range("I2").Value = 531.5
$ws.Range("K2").Value = 531.5
$ws.Range("M2").Value = -418.5
$ws.Range("H60").Value = 67527.336
$ws.Range("I60").Value = 70569.36
$ws.Range("J60").Value = 24939
$ws.Range("K60").Value = 70569.36
$ws.Range("L60").Value = 24939
$ws.Range("M60").Value = -69836.36
$ws.Range("N60").Value = -26405
$ws.Range("H88").Value = 4111.2856
$ws.Range("J88").Value = 4246.5
$ws.Range("L88").Value = 4246.5
$ws.Range("N88").Value = -5058.5
$ws.Range("H91").Value = 4111.2856
$ws.Range("J91").Value = 4246.5
$ws.Range("L91").Value = 4246.5
$ws.Range("N91").Value = -7054.5
$ws.Range("H116").Value = 531.5
$ws.Range("I116").Value = 531.5
$ws.Range("K116").Value = 531.5
$ws.Range("M116").Value = 1762.5
$ws.Range("H122").Value = 2342.9092
$ws.Range("I122").Value = 2039.05
$ws.Range("J122").Value = 5381.5
$ws.Range("K122").Value = 6117.15
$ws.Range("L122").Value = 16144.5
$ws.Range("M122").Value = -3667.15
$ws.Range("N122").Value = -21044.5
$ws.Range("H132").Value = 1727148.8
$ws.Range("I132").Value = 2939.8215
$ws.Range("K132").Value = 8819.4645
$ws.Range("M132").Value = -6289.4645

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 531.5
$ws.Range("I3").Value = 531.5
$ws.Range("K3").Value = 531.5
$ws.Range("M3").Value = -417.5
$ws.Range("H134").Value = 2859291.2
$ws.Range("I134").Value = 2083.2334
$ws.Range("K134").Value = 6249.7002
$ws.Range("M134").Value = -3714.7002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25643590
$ws.Range("I31").Value = 38463960
$ws.Range("J31").Value = 2849.5386
$ws.Range("K31").Value = 38463960
$ws.Range("L31").Value = 2849.5386
$ws.Range("M31").Value = -38463665
$ws.Range("N31").Value = -3439.5386
$ws.Range("H34").Value = 25643590
$ws.Range("I34").Value = 38463960
$ws.Range("J34").Value = 2849.5386
$ws.Range("K34").Value = 38463960
$ws.Range("L34").Value = 2849.5386
$ws.Range("M34").Value = -38463758
$ws.Range("N34").Value = -3253.5386
$ws.Range("H105").Value = 1355.6
$ws.Range("J105").Value = 1666.6666
$ws.Range("L105").Value = 1666.6666
$ws.Range("N105").Value = -5160.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1339.0667
$ws.Range("I12").Value = 138
$ws.Range("K12").Value = 414
$ws.Range("M12").Value = -241
$ws.Range("H16").Value = 21666
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 21666
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 64998
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -65344
$ws.Range("H86").Value = 363.75
$ws.Range("I86").Value = 344.2857
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 1032.8571
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = 153.1428999999998
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 363.75
$ws.Range("I89").Value = 344.2857
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 3098.5713
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = 2829.4287
$ws.Range("N89").Value = -16356
$ws.Range("H109").Value = 15339.444
$ws.Range("I109").Value = 1181.25
$ws.Range("K109").Value = 3543.75
$ws.Range("M109").Value = -2503.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 41166.668
$ws.Range("J26").Value = 41166.668
$ws.Range("L26").Value = 41166.668
$ws.Range("N26").Value = -41726.668
$ws.Range("H50").Value = 41166.668
$ws.Range("J50").Value = 41166.668
$ws.Range("L50").Value = 41166.668
$ws.Range("N50").Value = -42162.668
$ws.Range("H53").Value = 49998.5
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 49998.5
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 49998.5
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -51260.5
$ws.Range("H102").Value = 4000
$ws.Range("I102").Value = 4000
$ws.Range("K102").Value = 4000
$ws.Range("M102").Value = -2378

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8224.267
$ws.Range("J7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("N7").Value = -15224
$ws.Range("H22").Value = 9077.462
$ws.Range("I22").Value = 11472.7
$ws.Range("J22").Value = 1093.3334
$ws.Range("K22").Value = 11472.7
$ws.Range("L22").Value = 1093.3334
$ws.Range("M22").Value = -11177.7
$ws.Range("N22").Value = -1683.3334
$ws.Range("H23").Value = 28998.25
$ws.Range("J23").Value = 28998.25
$ws.Range("L23").Value = 28998.25
$ws.Range("N23").Value = -29458.25
$ws.Range("H27").Value = 9077.462
$ws.Range("I27").Value = 11472.7
$ws.Range("J27").Value = 1093.3334
$ws.Range("K27").Value = 11472.7
$ws.Range("L27").Value = 1093.3334
$ws.Range("M27").Value = -11365.7
$ws.Range("N27").Value = -1307.3334
$ws.Range("H40").Value = 6598.2
$ws.Range("I40").Value = 6748
$ws.Range("K40").Value = 6748
$ws.Range("M40").Value = -6612
$ws.Range("H46").Value = 2500.6667
$ws.Range("J46").Value = 2500.6667
$ws.Range("L46").Value = 2500.6667
$ws.Range("N46").Value = -2876.6667
$ws.Range("H56").Value = 33872.43
$ws.Range("I56").Value = 25410.2
$ws.Range("J56").Value = 55028
$ws.Range("K56").Value = 25410.2
$ws.Range("L56").Value = 55028
$ws.Range("M56").Value = -24719.2
$ws.Range("N56").Value = -56410
$ws.Range("H58").Value = 49991
$ws.Range("J58").Value = 49999
$ws.Range("L58").Value = 49999
$ws.Range("N58").Value = -50519
$ws.Range("H68").Value = 2453108.2
$ws.Range("I68").Value = 2978246
$ws.Range("J68").Value = 2466
$ws.Range("K68").Value = 2978246
$ws.Range("L68").Value = 2466
$ws.Range("M68").Value = -2977497
$ws.Range("N68").Value = -3964
$ws.Range("H71").Value = 2453108.2
$ws.Range("I71").Value = 2978246
$ws.Range("J71").Value = 2466
$ws.Range("K71").Value = 14891230
$ws.Range("L71").Value = 12330
$ws.Range("M71").Value = -14887486
$ws.Range("N71").Value = -19818
$ws.Range("H82").Value = 3938.3635
$ws.Range("J82").Value = 4510.5
$ws.Range("L82").Value = 4510.5
$ws.Range("N82").Value = -5232.5
$ws.Range("H85").Value = 3938.3635
$ws.Range("J85").Value = 4510.5
$ws.Range("L85").Value = 4510.5
$ws.Range("N85").Value = -7006.5
$ws.Range("H126").Value = 8224.267
$ws.Range("J126").Value = 15000
$ws.Range("L126").Value = 45000
$ws.Range("N126").Value = -49940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 20000
$ws.Range("I58").Value = 20000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 20000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -19692
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 29995
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H81").Value = 1995
$ws.Range("I81").Value = 1990
$ws.Range("K81").Value = 3980
$ws.Range("M81").Value = -2919
$ws.Range("H84").Value = 1995
$ws.Range("I84").Value = 1990
$ws.Range("K84").Value = 19900
$ws.Range("M84").Value = -14596
$ws.Range("H107").Value = 2137.65
$ws.Range("I107").Value = 1183.84
$ws.Range("J107").Value = 3727.3333
$ws.Range("K107").Value = 3551.52
$ws.Range("L107").Value = 11181.9999
$ws.Range("M107").Value = -1631.52
$ws.Range("N107").Value = -15021.9999
$ws.Range("H122").Value = 2117.8333
$ws.Range("I122").Value = 1843.1177
$ws.Range("J122").Value = 6788
$ws.Range("K122").Value = 5529.3531
$ws.Range("L122").Value = 20364
$ws.Range("M122").Value = -3079.3531
$ws.Range("N122").Value = -25264
$ws.Range("H136").Value = 169289.44
$ws.Range("I136").Value = 2695.9656
$ws.Range("K136").Value = 8087.8968
$ws.Range("M136").Value = -5537.8968
